# Atualizando o arquivo XLSX
# Refreshes Betfair back/lay odds figures (columns F:AO) for the
# matches listed in Jogos_do_Dia_Betfair_Back_Lay_2025-11-26.xlsx.
# Only numeric odds cells change; League/Date/Time/Home/Away stay put.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 6).Value = 5.3
$ws.Cells.Item(2, 7).Value = 7.2
$ws.Cells.Item(2, 8).Value = 1.52
$ws.Cells.Item(2, 9).Value = 1.73
$ws.Cells.Item(2, 10).Value = 3.6
$ws.Cells.Item(2, 12).Value = 1.27
$ws.Cells.Item(2, 13).Value = 1.04
$ws.Cells.Item(2, 14).Value = 4
$ws.Cells.Item(2, 16).Value = 2.16
$ws.Cells.Item(2, 19).Value = 2.48
$ws.Cells.Item(2, 20).Value = 1.78
$ws.Cells.Item(2, 22).Value = 2.38
$ws.Cells.Item(2, 24).Value = 21
$ws.Cells.Item(2, 25).Value = 12
$ws.Cells.Item(2, 26).Value = 12.5
$ws.Cells.Item(2, 27).Value = 17.5
$ws.Cells.Item(2, 28).Value = 27
$ws.Cells.Item(2, 29).Value = 12.5
$ws.Cells.Item(2, 30).Value = 12.5
$ws.Cells.Item(2, 31).Value = 19
$ws.Cells.Item(2, 32).Value = 70
$ws.Cells.Item(2, 33).Value = 29
$ws.Cells.Item(2, 34).Value = 25
$ws.Cells.Item(2, 35).Value = 40
$ws.Cells.Item(2, 41).Value = 9.4

# Row 4
$ws.Cells.Item(4, 6).Value = 1.35
$ws.Cells.Item(4, 7).Value = 1.36
$ws.Cells.Item(4, 14).Value = 5.8
$ws.Cells.Item(4, 16).Value = 2.62
$ws.Cells.Item(4, 17).Value = 1.58
$ws.Cells.Item(4, 18).Value = 1.66
$ws.Cells.Item(4, 19).Value = 2.46
$ws.Cells.Item(4, 20).Value = 1.92
$ws.Cells.Item(4, 21).Value = 2.02
$ws.Cells.Item(4, 23).Value = 3.75
$ws.Cells.Item(4, 25).Value = 40
$ws.Cells.Item(4, 28).Value = 10.5
$ws.Cells.Item(4, 31).Value = 150
$ws.Cells.Item(4, 34).Value = 25

# Row 5
$ws.Cells.Item(5, 8).Value = 1.73
$ws.Cells.Item(5, 9).Value = 1.74
$ws.Cells.Item(5, 21).Value = 2.26
$ws.Cells.Item(5, 22).Value = 2.34
$ws.Cells.Item(5, 24).Value = 21
$ws.Cells.Item(5, 28).Value = 22
$ws.Cells.Item(5, 30).Value = 10
$ws.Cells.Item(5, 32).Value = 42
$ws.Cells.Item(5, 36).Value = 130
$ws.Cells.Item(5, 37).Value = 60
$ws.Cells.Item(5, 39).Value = 85
$ws.Cells.Item(5, 41).Value = 8

# Row 6
$ws.Cells.Item(6, 10).Value = 1.03
$ws.Cells.Item(6, 22).Value = 1.48

# Row 7
$ws.Cells.Item(7, 6).Value = 2.1
$ws.Cells.Item(7, 7).Value = 2.22
$ws.Cells.Item(7, 9).Value = 3.8
$ws.Cells.Item(7, 12).Value = 1.3
$ws.Cells.Item(7, 20).Value = 1.56
$ws.Cells.Item(7, 23).Value = 1.82
$ws.Cells.Item(7, 25).Value = 24
$ws.Cells.Item(7, 26).Value = 36
$ws.Cells.Item(7, 28).Value = 16.5
$ws.Cells.Item(7, 29).Value = 11.5
$ws.Cells.Item(7, 31).Value = 44
$ws.Cells.Item(7, 33).Value = 13.5
$ws.Cells.Item(7, 35).Value = 48
$ws.Cells.Item(7, 37).Value = 25
$ws.Cells.Item(7, 38).Value = 34
$ws.Cells.Item(7, 40).Value = 13
$ws.Cells.Item(7, 41).Value = 32

# Row 8
$ws.Cells.Item(8, 12).Value = 1.38
$ws.Cells.Item(8, 16).Value = 2.08
$ws.Cells.Item(8, 17).Value = 1.9
$ws.Cells.Item(8, 19).Value = 3.25
$ws.Cells.Item(8, 24).Value = 16
$ws.Cells.Item(8, 29).Value = 7.8
$ws.Cells.Item(8, 34).Value = 16.5
$ws.Cells.Item(8, 41).Value = 30

# Row 9
$ws.Cells.Item(9, 8).Value = 9.6
$ws.Cells.Item(9, 11).Value = 6.8
$ws.Cells.Item(9, 16).Value = 3.6
$ws.Cells.Item(9, 21).Value = 2.46
$ws.Cells.Item(9, 23).Value = 3.9
$ws.Cells.Item(9, 25).Value = 55
$ws.Cells.Item(9, 27).Value = 290
$ws.Cells.Item(9, 28).Value = 16.5

# Row 10
$ws.Cells.Item(10, 8).Value = 9.4
$ws.Cells.Item(10, 9).Value = 9.800000000000001
$ws.Cells.Item(10, 10).Value = 5.6
$ws.Cells.Item(10, 11).Value = 5.7
$ws.Cells.Item(10, 16).Value = 2.7
$ws.Cells.Item(10, 17).Value = 1.55
$ws.Cells.Item(10, 21).Value = 2.12
$ws.Cells.Item(10, 22).Value = 1.11
$ws.Cells.Item(10, 24).Value = 26
$ws.Cells.Item(10, 30).Value = 32
$ws.Cells.Item(10, 32).Value = 9.4
$ws.Cells.Item(10, 41).Value = 110

# Row 11
$ws.Cells.Item(11, 14).Value = 5.7
$ws.Cells.Item(11, 18).Value = 1.65
$ws.Cells.Item(11, 20).Value = 1.64
$ws.Cells.Item(11, 23).Value = 2.42
$ws.Cells.Item(11, 25).Value = 26
$ws.Cells.Item(11, 39).Value = 75

# Row 12
$ws.Cells.Item(12, 12).Value = 1.35
$ws.Cells.Item(12, 19).Value = 2.94
$ws.Cells.Item(12, 20).Value = 1.65
$ws.Cells.Item(12, 21).Value = 2.44
$ws.Cells.Item(12, 27).Value = 55
$ws.Cells.Item(12, 33).Value = 10.5
$ws.Cells.Item(12, 36).Value = 29

# Row 13
$ws.Cells.Item(13, 6).Value = 7.6
$ws.Cells.Item(13, 7).Value = 7.8
$ws.Cells.Item(13, 8).Value = 1.47
$ws.Cells.Item(13, 9).Value = 1.48
$ws.Cells.Item(13, 17).Value = 1.56
$ws.Cells.Item(13, 19).Value = 2.34
$ws.Cells.Item(13, 20).Value = 1.75
$ws.Cells.Item(13, 21).Value = 2.24
$ws.Cells.Item(13, 22).Value = 3.1
$ws.Cells.Item(13, 23).Value = 1.15
$ws.Cells.Item(13, 25).Value = 12
$ws.Cells.Item(13, 27).Value = 14
$ws.Cells.Item(13, 29).Value = 12
$ws.Cells.Item(13, 31).Value = 13.5
$ws.Cells.Item(13, 34).Value = 20
$ws.Cells.Item(13, 36).Value = 250
$ws.Cells.Item(13, 38).Value = 75
$ws.Cells.Item(13, 39).Value = 85
$ws.Cells.Item(13, 40).Value = 80
$ws.Cells.Item(13, 41).Value = 5.1

# Row 14
$ws.Cells.Item(14, 6).Value = 2.72
$ws.Cells.Item(14, 7).Value = 2.76
$ws.Cells.Item(14, 8).Value = 2.68
$ws.Cells.Item(14, 9).Value = 2.72
$ws.Cells.Item(14, 22).Value = 1.58
$ws.Cells.Item(14, 41).Value = 16.5

# Row 16
$ws.Cells.Item(16, 6).Value = 1.95
$ws.Cells.Item(16, 7).Value = 2
$ws.Cells.Item(16, 8).Value = 4.1
$ws.Cells.Item(16, 10).Value = 3.75
$ws.Cells.Item(16, 16).Value = 2.08
$ws.Cells.Item(16, 17).Value = 1.75
$ws.Cells.Item(16, 18).Value = 1.43
$ws.Cells.Item(16, 20).Value = 1.73
$ws.Cells.Item(16, 21).Value = 2.18
$ws.Cells.Item(16, 23).Value = 2
$ws.Cells.Item(16, 24).Value = 21
$ws.Cells.Item(16, 27).Value = 110
$ws.Cells.Item(16, 32).Value = 13
$ws.Cells.Item(16, 34).Value = 18
$ws.Cells.Item(16, 40).Value = 12.5

# Row 18
$ws.Cells.Item(18, 6).Value = 1.79
$ws.Cells.Item(18, 7).Value = 1.94
$ws.Cells.Item(18, 8).Value = 4.9
$ws.Cells.Item(18, 10).Value = 3.4
$ws.Cells.Item(18, 12).Value = 1.48
$ws.Cells.Item(18, 13).Value = 1.08
$ws.Cells.Item(18, 16).Value = 1.67
$ws.Cells.Item(18, 17).Value = 2.18
$ws.Cells.Item(18, 19).Value = 3.75
$ws.Cells.Item(18, 20).Value = 2
$ws.Cells.Item(18, 21).Value = 1.8
$ws.Cells.Item(18, 23).Value = 2.06
$ws.Cells.Item(18, 35).Value = 120

# Row 19
$ws.Cells.Item(19, 6).Value = 2.12
$ws.Cells.Item(19, 7).Value = 2.28
$ws.Cells.Item(19, 8).Value = 4.5
$ws.Cells.Item(19, 9).Value = 5.1
$ws.Cells.Item(19, 10).Value = 2.92
$ws.Cells.Item(19, 11).Value = 3.1
$ws.Cells.Item(19, 12).Value = 1.62
$ws.Cells.Item(19, 13).Value = 1.12
$ws.Cells.Item(19, 16).Value = 1.43
$ws.Cells.Item(19, 17).Value = 2.88
$ws.Cells.Item(19, 18).Value = 1.15
$ws.Cells.Item(19, 19).Value = 6.2
$ws.Cells.Item(19, 20).Value = 2.28
$ws.Cells.Item(19, 21).Value = 1.63
$ws.Cells.Item(19, 22).Value = 1.24
$ws.Cells.Item(19, 23).Value = 1.78
$ws.Cells.Item(19, 26).Value = 980
$ws.Cells.Item(19, 28).Value = 7.4
$ws.Cells.Item(19, 29).Value = 8.6
$ws.Cells.Item(19, 32).Value = 12.5
$ws.Cells.Item(19, 33).Value = 13
$ws.Cells.Item(19, 35).Value = 160
$ws.Cells.Item(19, 38).Value = 80
$ws.Cells.Item(19, 39).Value = 320
$ws.Cells.Item(19, 40).Value = 1000

# Row 20
$ws.Cells.Item(20, 6).Value = 3.55
$ws.Cells.Item(20, 7).Value = 4.2
$ws.Cells.Item(20, 8).Value = 2.08
$ws.Cells.Item(20, 9).Value = 2.26
$ws.Cells.Item(20, 10).Value = 3.5
$ws.Cells.Item(20, 12).Value = 1.32
$ws.Cells.Item(20, 14).Value = 3.85
$ws.Cells.Item(20, 15).Value = 1.3
$ws.Cells.Item(20, 16).Value = 1.98
$ws.Cells.Item(20, 18).Value = 1.38
$ws.Cells.Item(20, 19).Value = 3.25
$ws.Cells.Item(20, 22).Value = 1.79
$ws.Cells.Item(20, 23).Value = 1.33
$ws.Cells.Item(20, 26).Value = 14.5
$ws.Cells.Item(20, 28).Value = 1000
$ws.Cells.Item(20, 33).Value = 1000
$ws.Cells.Item(20, 41).Value = 16.5

# Row 21
$ws.Cells.Item(21, 7).Value = 1.76
$ws.Cells.Item(21, 8).Value = 5.1
$ws.Cells.Item(21, 9).Value = 7
$ws.Cells.Item(21, 10).Value = 3.6
$ws.Cells.Item(21, 11).Value = 5.1
$ws.Cells.Item(21, 14).Value = 4.2
$ws.Cells.Item(21, 16).Value = 2.06
$ws.Cells.Item(21, 17).Value = 1.74
$ws.Cells.Item(21, 18).Value = 1.42
$ws.Cells.Item(21, 19).Value = 2.92
$ws.Cells.Item(21, 20).Value = 1.76
$ws.Cells.Item(21, 21).Value = 2.02
$ws.Cells.Item(21, 22).Value = 1.18
$ws.Cells.Item(21, 23).Value = 2.3
$ws.Cells.Item(21, 33).Value = 12.5
$ws.Cells.Item(21, 39).Value = 130
